$wb = $excel.ActiveWorkbook

# Old and new "built on" timestamps embedded in the version strings.
$oldStamp = "January 30 2026 16.19.47 EST"
$newStamp = "February 02 2026 12.49.33 EST"

$aboutSheet = $wb.Worksheets.Item("About")
$dataSheet  = $wb.Worksheets.Item("Boundaries and methane sources")

# "About" sheet: version banner (A2) and recommended citation (A6).
$aboutSheet.Range("A2").Value = "Version: mines - January 30 (built on $newStamp)"
$aboutSheet.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Huoerxinhe Coal Mine, China, M1140, version 'mines - January 30 (built on $newStamp)'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# "Boundaries and methane sources" sheet: build_version column (S), rows 2-8.
for ($r = 2; $r -le 8; $r++) {
    $cell = $dataSheet.Cells.Item($r, 19)
    if ($cell.Value() -eq "mines - January 30 (built on $oldStamp)") {
        $cell.Value = "mines - January 30 (built on $newStamp)"
    }
}
